# Redid estimations of hysplits so it only extends back two days
# Updates the D1 (col C) / D2 (col D) compass-direction estimates for
# several rows, marks the revised rows (54:57) with red text, moves the
# active selection/scroll position, and sets the page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Day-1 (column C) / Day-2 (column D) direction updates ---------------

# Rows 2-5: Day1 N -> S
$ws.Range("C2:C5").Value = "S"

# Rows 10-13: Day1 W -> S, Day2 newly recorded as W
$ws.Range("C10:C13").Value = "S"
$ws.Range("D10:D13").Value = "W"

# Rows 26-29: Day1 E -> N, Day2 N -> E
$ws.Range("C26:C29").Value = "N"
$ws.Range("D26:D29").Value = "E"

# Rows 54-56: Day1 W -> S, Day2 -> W
$ws.Range("C54:C56").Value = "S"
$ws.Range("D54:D56").Value = "W"

# Row 57: Day1 E -> S (Day2 stays W)
$ws.Range("C57").Value = "S"

# Rows 58-61: Day2 newly recorded as N (Day1 W unchanged)
$ws.Range("D58:D61").Value = "N"

# Rows 62-65: Day1 E -> N
$ws.Range("C62:C65").Value = "N"

# Rows 78-81: Day2 newly recorded as S (Day1 E unchanged)
$ws.Range("D78:D81").Value = "S"

# Rows 85-89: Day1 E -> S
$ws.Range("C85:C89").Value = "S"

# Rows 90-91: Day1 W -> S
$ws.Range("C90:C91").Value = "S"

# --- Highlight the re-estimated rows 54:57 in red -------------------------
$ws.Range("C54:C57").Font.Color = 255

# --- Page setup ------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Update scroll position / selection ------------------------------------
$excel.ActiveWindow.ScrollRow = 74
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C81").Select()

"done"
